# Actualización automática 2025-05-29 10:10:07
#
# A new salesperson row ("OFICINA-CATAECSA" / "MACHARE BARCO LISSETTE
# STEFANIA") is inserted at row 240 on both worksheets, pushing the
# existing rows 240-260 down to 241-261 (and the trailing totals row
# from 261 down to 262). The new row is populated with zero amounts.
# Finally the "x de 259" counter labels in the grand-total row of the
# first sheet are refreshed to "x de 260" to reflect the extra record.

$wb = $excel.ActiveWorkbook

$sheetNames = @("VENTAS POR GRUPO", "VENTA MENSUAL")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new blank row before row 240 - this shifts every row at
    # or below 240 down by one and copies formatting from the row above.
    $ws.Rows.Item(240).Insert()

    # Determine how many numeric data columns this sheet has (N on the
    # "VENTAS POR GRUPO" sheet, F on "VENTA MENSUAL").
    if ($sheetName -eq "VENTAS POR GRUPO") {
        $lastCol = 14
    } else {
        $lastCol = 6
    }

    $ws.Cells.Item(240, 1).Value = "OFICINA-CATAECSA"
    $ws.Cells.Item(240, 2).Value = "MACHARE BARCO LISSETTE STEFANIA"

    for ($col = 3; $col -le $lastCol; $col++) {
        $ws.Cells.Item(240, $col).Value = 0
    }
}

# Refresh the "x de 259" -> "x de 260" counters in the grand-total row
# of the "VENTAS POR GRUPO" sheet (the totals row shifted from 261 to
# 262 because of the inserted row above).
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
for ($col = 3; $col -le 14; $col++) {
    $cell = $wsGrupo.Cells.Item(262, $col)
    $text = $cell.Value()
    $cell.Value = $text.Replace("de 259", "de 260")
}
